$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(48).Insert()

$ws.Cells.Item(48, 1).Value = 10
$ws.Cells.Item(48, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(48, 3).Value = "La Araucanía"
$ws.Cells.Item(48, 4).Value = 45272
$ws.Cells.Item(48, 5).Value = 9
$ws.Cells.Item(48, 6).Value = "Fruta"
$ws.Cells.Item(48, 7).Value = 100108
$ws.Cells.Item(48, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(48, 9).Value = 100108007
$ws.Cells.Item(48, 10).Value = "Coco"
$ws.Cells.Item(48, 11).Value = "Sin especificar"
$ws.Cells.Item(48, 12).Value = "Primera"
$ws.Cells.Item(48, 13).Value = 40
$ws.Cells.Item(48, 14).Value = 25000
$ws.Cells.Item(48, 15).Value = 25000
$ws.Cells.Item(48, 16).Value = 25000
$ws.Cells.Item(48, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(48, 18).Value = "Perú"
$ws.Cells.Item(48, 19).Value = 1250
$ws.Cells.Item(48, 20).Value = 20
